$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "AI Co-Op"
$ws.Range("B2").Value = "ArtiFlex Manufacturing"
$ws.Range("C2").Value = "Wooster, OH, US USA"
$ws.Range("D2").Value = 24.4
$ws.Range("E2").Value = "Generative AI, LangChain, RAG, LLaMA, Azure ML, S3, Glue, Athena, Redshift, BigQuery"
$ws.Range("G2").Value = "https://www.indeed.com/viewjob?jk=9bd105687c5aa984"

# Row 3
$ws.Range("A3").Value = "Senior Systems Engineer, UDS Data Management"
$ws.Range("B3").Value = "Dell Technologies"
$ws.Range("C3").Value = "Remote, US USA"
$ws.Range("D3").Value = 23.3
$ws.Range("E3").Value = "Data Scientist, LangChain, RAG, LLaMA, TensorFlow, PyTorch, Redshift, BigQuery, Synapse, Data Lake"
$ws.Range("G3").Value = "https://www.indeed.com/viewjob?jk=8b555f123481de91"

# Row 4
$ws.Range("A4").Value = "ML Ops Engineer"
$ws.Range("B4").Value = "Hitachi Digital Services"
$ws.Range("C4").Value = "Reading, PA, US USA"
$ws.Range("D4").Value = 21.1
$ws.Range("E4").Value = "Data Scientist, RAG, TensorFlow, S3, MLflow, Docker, Kubernetes, CI/CD, Jenkins, GitHub Actions"
$ws.Range("G4").Value = "https://www.indeed.com/viewjob?jk=5df2e993fb881155"

# Row 5
$ws.Range("A5").Value = "Data Scientist- Payments Optimization"
$ws.Range("B5").Value = "Worldpay"
$ws.Range("C5").Value = "Cincinnati, OH, US USA"
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = "Data Scientist, TensorFlow, CI/CD, Git, Snowflake, Databricks, PySpark, Hadoop, Tableau, Power BI"
$ws.Range("G5").Value = "https://www.indeed.com/viewjob?jk=68e605c94d2c096e"

# Row 6
$ws.Range("A6").Value = "Software Development Engineer"
$ws.Range("B6").Value = "Expedia Group"
$ws.Range("C6").Value = "San Jose, CA, US USA"
$ws.Range("D6").Value = 15.6
$ws.Range("E6").Value = "RAG, S3, EC2, Kubernetes, CI/CD, Jenkins, GitHub Actions, Terraform, Git, Python"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "2026-02-23"
$ws.Range("G6").Value = "https://www.indeed.com/viewjob?jk=f6610ed807a1c6f6"

# Row 7
$ws.Range("A7").Value = "Senior Observability Platform Engineer"
$ws.Range("B7").Value = "Klaviyo"
$ws.Range("C7").Value = "Boston, MA, US USA"
$ws.Range("D7").Value = 14.4
$ws.Range("E7").Value = "RAG, Cortex, Kubernetes, CI/CD, Terraform, Git, Kafka, MySQL, Python, SQL"
$ws.Range("G7").Value = "https://www.indeed.com/viewjob?jk=6db8b13f873a0ab3"

# Row 8
$ws.Range("A8").Value = "AWS AI & DevOps Intern"
$ws.Range("B8").Value = "Network Distribution"
$ws.Range("C8").Value = "Schaumburg, IL, US USA"
$ws.Range("D8").Value = 14.4
$ws.Range("E8").Value = "Generative AI, Copilot, TensorFlow, PyTorch, S3, EC2, Glue, CI/CD, Terraform, Git"
$ws.Range("G8").Value = "https://www.indeed.com/viewjob?jk=f896c6da0dfabb9a"

# Row 9
$ws.Range("A9").Value = "Information Technology - BI Data Architect"
$ws.Range("B9").Value = "TCC Verizon Authorized Retailer"
$ws.Range("C9").Value = "Fishers, IN, US USA"
$ws.Range("D9").Value = 13.3
$ws.Range("E9").Value = "Data Scientist, Copilot, Synapse, Git, Databricks, PySpark, Power BI, Python, SQL, R"
$ws.Range("G9").Value = "https://www.indeed.com/viewjob?jk=a149a432f720bb1c"

# Row 10
$ws.Range("A10").Value = "Sr Data Architect"
$ws.Range("B10").Value = "Bank of America"
$ws.Range("C10").Value = "Charlotte, NC, US USA"
$ws.Range("E10").Value = "AI Engineer, Data Scientist, RAG, Synapse, Data Lake, Snowflake, Databricks, Python, SQL, R"
$ws.Range("G10").Value = "https://www.indeed.com/viewjob?jk=3c087e00bef6a423"

# Row 11
$ws.Range("A11").Value = "Software Engineer, Financial Systems"
$ws.Range("B11").Value = "Opendoor"
$ws.Range("C11").Value = "Seattle, WA, US USA"
$ws.Range("E11").Value = "Docker, CI/CD, GitHub Actions, Terraform, Git, Snowflake, Quicksight, Python, SQL, R"
$ws.Range("G11").Value = "https://www.indeed.com/viewjob?jk=c304968f71791ba5"

# Row 12
$ws.Range("A12").Value = "Senior Software Engineer (AI Enablement)"
$ws.Range("B12").Value = "Redwood Logistics"
$ws.Range("C12").Value = "Remote, US USA"
$ws.Range("E12").Value = "RAG, AKS, CI/CD, Git, Snowflake, Kafka, SQL, R, Java, Scala"
$ws.Range("G12").Value = "https://www.indeed.com/viewjob?jk=9f0a71f110d4bcaf"

# Row 13
$ws.Range("A13").Value = "Software Engineer"
$ws.Range("B13").Value = "NetSPI"
$ws.Range("C13").Value = "Minneapolis, MN, US USA"
$ws.Range("E13").Value = "S3, Docker, GitHub Actions, Git, PostgreSQL, Python, SQL, R, Java, Scala"
$ws.Range("G13").Value = "https://www.indeed.com/viewjob?jk=64c43e0c9ccd4fc4"

# Row 14
$ws.Range("A14").Value = "Senior Gen AI Application Developer - LangGraph, Bedrock, AWS (ONSITE)"
$ws.Range("B14").Value = "Cognizant Technology Solutions"
$ws.Range("C14").Value = "Detroit, MI, US USA"
$ws.Range("E14").Value = "AI Engineer, Data Scientist, Generative AI, LangChain, RAG, FAISS, Pinecone, CI/CD, R, Scala"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "2026-02-18"
$ws.Range("G14").Value = "https://www.indeed.com/viewjob?jk=472f11358fc2e3f7"

# Row 15
$ws.Range("A15").Value = "Senior Data Scientist"
$ws.Range("B15").Value = "Guild Mortgage Company LLC"
$ws.Range("C15").Value = "US USA"
$ws.Range("E15").Value = "Data Scientist, RAG, TensorFlow, Data Lake, Hadoop, Tableau, Python, SQL, R, Scala"
$ws.Range("G15").Value = "https://www.indeed.com/viewjob?jk=b3b9c468aaf3a5d0"

# Row 16
$ws.Range("A16").Value = "Associate Data Scientist, New College Grad - 2026"
$ws.Range("B16").Value = "Visa"
$ws.Range("C16").Value = "Foster City, CA, US USA"
$ws.Range("E16").Value = "Data Scientist, Generative AI, RAG, Prompt Engineering, Git, Tableau, Power BI, Python, SQL, R"
$ws.Range("G16").Value = "https://www.indeed.com/viewjob?jk=45d3814a2a03d341"

# Row 17
$ws.Range("A17").Value = "Data Engineer (Python/Spark)"
$ws.Range("B17").Value = "Take-Two Interactive Software, Inc."
$ws.Range("C17").Value = "Austin, TX, US USA"
$ws.Range("D17").Value = 11.1
$ws.Range("E17").Value = "RAG, S3, Docker, Jenkins, Git, Python, SQL, R, Java, Scala"
$ws.Range("G17").Value = "https://www.indeed.com/viewjob?jk=98762eeaa2f684be"

# Row 18
$ws.Range("A18").Value = "Sr Software Engineer - Remote"
$ws.Range("B18").Value = "Optum"
$ws.Range("C18").Value = "Basking Ridge, NJ, US USA"
$ws.Range("E18").Value = "Docker, Kubernetes, Jenkins, PySpark, Hadoop, Python, SQL, R, Optimization"
$ws.Range("G18").Value = "https://www.indeed.com/viewjob?jk=d1aa56bb760b5544"
